$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-mark cells whose new numeric-looking text must stay as TEXT (not be
# auto-converted to a Number), matching the original inlineStr/text storage.
$textFormatRows = @(5, 6, 9, 10, 12, 15, 16, 19, 20, 21, 22, 24, 25, 29, 32, 33, 34, 35, 36, 37, 38, 39, 40, 42, 43, 44, 45, 46, 49, 50, 51)
foreach ($r in $textFormatRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = '57.753.83'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '3.113.52'
$ws.Range("E3").Value = '  +1.40%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '524.88'
$ws.Range("E5").Value = '  +1.76%  '
$ws.Range("D6").Value = '141.81'
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.113.26'
$ws.Range("E8").Value = '  +1.41%  '
$ws.Range("D9").Value = '0.436'
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").Value = '7.30'
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("D12").Value = '0.383'
$ws.Range("E12").Value = '  +2.56%  '
$ws.Range("D13").Value = '3.641.72'
$ws.Range("E13").Value = '  +0.97%  '
$ws.Range("E14").Value = '  +1.88%  '
$ws.Range("D15").Value = '26.25'
$ws.Range("E15").Value = '  +2.71%  '
$ws.Range("D16").Value = '0.0000165'
$ws.Range("E16").Value = '  +1.39%  '
$ws.Range("D17").Value = '57.846.16'
$ws.Range("D18").Value = '3.107.61'
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").Value = '6.11'
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("D20").Value = '12.84'
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("D21").Value = '8.07'
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").Value = '338.16'
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '0.508'
$ws.Range("E24").Value = '  +1.46%  '
$ws.Range("D25").Value = '66.34'
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("D28").Value = '0.0₃0932'
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").Value = '6.58'
$ws.Range("E29").Value = '  +3.70%  '
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").Value = '1.87'
$ws.Range("E32").Value = '  +2.52%  '
$ws.Range("D33").Value = '1.21'
$ws.Range("E33").Value = '  +3.85%  '
$ws.Range("D34").Value = '20.94'
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("D35").Value = '154.47'
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").Value = '4.61'
$ws.Range("E36").Value = '  +3.48%  '
$ws.Range("D37").Value = '6.06'
$ws.Range("E37").Value = '  +3.49%  '
$ws.Range("D38").Value = '26.99'
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").Value = '1.30'
$ws.Range("E39").Value = '  +1.46%  '
$ws.Range("D40").Value = '0.0667'
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("D41").Value = '3.151.34'
$ws.Range("E41").Value = '  +0.96%  '
$ws.Range("D42").Value = '0.684'
$ws.Range("E42").Value = '  +3.85%  '
$ws.Range("D43").Value = '3.91'
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("D44").Value = '36.90'
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").Value = '1.47'
$ws.Range("E46").Value = '  +6.50%  '
$ws.Range("D47").Value = '2.279.42'
$ws.Range("E47").Value = '  +0.83%  '
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("D49").Value = '0.972'
$ws.Range("E49").Value = '  +5.05%  '
$ws.Range("D50").Value = '20.68'
$ws.Range("E50").Value = '  +3.53%  '
$ws.Range("D51").Value = '6.02'
$ws.Range("E51").Value = '  +3.04%  '

# Restore the default (Normal) style on those cells so the underlying
# number-format metadata matches the original workbook (text value only,
# no visible/style change).
foreach ($r in $textFormatRows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}
